$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 59, shifting existing rows 59:133 down to 60:134
$ws.Rows(59).Insert()

# Populate the newly inserted row 59 with the new record's data
$ws.Cells.Item(59, 1).Value = 4
$ws.Cells.Item(59, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value = "Los Lagos"
$ws.Cells.Item(59, 4).Value = (Get-Date -Year 2022 -Month 8 -Day 30 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(59, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(59, 5).Value = 10
$ws.Cells.Item(59, 6).Value = 100112022
$ws.Cells.Item(59, 7).Value = "Arveja Verde"
$ws.Cells.Item(59, 8).Value = "Perfection"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 70
$ws.Cells.Item(59, 11).Value = 45000
$ws.Cells.Item(59, 12).Value = 45000
$ws.Cells.Item(59, 13).Value = 45000
$ws.Cells.Item(59, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(59, 16).Value = 1800
$ws.Cells.Item(59, 17).Value = 25
$ws.Cells.Item(59, 18).Value = "Hortaliza"
